# Applies the data refresh described by the commit:
#   "Updated cryptos list on Wed Nov 13 04:30:46 UTC 2024 with GitHub
#    Actions"
#
# The sheet is a coinranking.com price/volume scrape: header in row 1,
# 50 coins in rows 2-51 (A=rank, B=coin name, C=link, D=price,
# E=1h volume %). This run refreshes the Price/Volume(1h) columns for
# (almost) every coin, and rows 47-48 swap contents: the coin that was
# "OKB" (row 47) is now "Stellar" and vice versa for row 48 (the
# underlying ranking data moved, not just the price).
#
# Every data cell in columns B:E is literal text (inlineStr) in the
# source workbook, even Price-column values that look like plain
# numbers ("204.15", "1.00", ...). Assigning such a string straight to
# `Range.Value` makes Excel's COM layer auto-coerce it into a *number*
# cell (e.g. "1.00" -> 1, losing the trailing zero) instead of leaving
# it as text, so for any replacement string that is itself a bare
# number we use the classic quote-prefix trick (a leading ' forces
# text entry the way typing it in the UI would) and then
# ClearFormats() so the transient "quote prefix" cell style left
# behind doesn't linger. Strings that aren't bare numbers (the "%"
# volume cells, and Price cells like "87.446.15" that contain more
# than one '.') already round-trip as text and are set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $Text) {
    $range = $ws.Range($CellRef)
    if ($Text.Trim() -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $Text
        $range.ClearFormats()
    } else {
        $range.Value = $Text
    }
}


# Row 2 (Bitcoin): D2, E2
Set-TextValue 'D2' '87.385.22'
Set-TextValue 'E2' '  -1.32%  '

# Row 3 (Ethereum): D3, E3
Set-TextValue 'D3' '3.215.94'
Set-TextValue 'E3' '  -3.63%  '

# Row 4 (TetherUSD): E4
Set-TextValue 'E4' '  +0.33%  '

# Row 5 (Solana): D5, E5
Set-TextValue 'D5' '203.83'
Set-TextValue 'E5' '  -7.45%  '

# Row 6 (BNB): D6, E6
Set-TextValue 'D6' '609.38'
Set-TextValue 'E6' '  -6.35%  '

# Row 7 (Dogecoin): D7, E7
Set-TextValue 'D7' '0.375'
Set-TextValue 'E7' '  -0.83%  '

# Row 8 (XRP): D8, E8
Set-TextValue 'D8' '0.671'
Set-TextValue 'E8' '  +10.66%  '

# Row 9 (USDC): D9, E9
Set-TextValue 'D9' '1.00'
Set-TextValue 'E9' '  +0.07%  '

# Row 10 (LidoStakedEther): D10, E10
Set-TextValue 'D10' '3.207.95'
Set-TextValue 'E10' '  -3.70%  '

# Row 11 (Cardano): E11
Set-TextValue 'E11' '  -7.94%  '

# Row 12 (TRON): E12
Set-TextValue 'E12' '  +6.56%  '

# Row 13 (ShibaInu): D13, E13
Set-TextValue 'D13' '0.0000245'
Set-TextValue 'E13' '  -13.07%  '

# Row 14 (WrappedliquidstakedEther2.0): D14, E14
Set-TextValue 'D14' '3.819.45'
Set-TextValue 'E14' '  -3.23%  '

# Row 15 (Toncoin): E15
Set-TextValue 'E15' '  -3.61%  '

# Row 16 (Avalanche): D16, E16
Set-TextValue 'D16' '32.63'
Set-TextValue 'E16' '  -8.06%  '

# Row 17 (WrappedBTC): D17, E17
Set-TextValue 'D17' '87.336.61'
Set-TextValue 'E17' '  -0.94%  '

# Row 18 (WrappedEther): D18, E18
Set-TextValue 'D18' '3.229.13'
Set-TextValue 'E18' '  -2.75%  '

# Row 19 (Chainlink): E19
Set-TextValue 'E19' '  -8.09%  '

# Row 20 (SuiNetwork): D20, E20
Set-TextValue 'D20' '2.94'
Set-TextValue 'E20' '  -6.52%  '

# Row 21 (BitcoinCash): D21, E21
Set-TextValue 'D21' '421.40'
Set-TextValue 'E21' '  -7.85%  '

# Row 22 (Uniswap): E22
Set-TextValue 'E22' '  -11.59%  '

# Row 23 (Polkadot): D23, E23
Set-TextValue 'D23' '5.11'
Set-TextValue 'E23' '  -6.91%  '

# Row 24 (NEARProtocol): D24, E24
Set-TextValue 'D24' '5.25'
Set-TextValue 'E24' '  -6.20%  '

# Row 25 (Aptos): E25
Set-TextValue 'E25' '  -8.37%  '

# Row 26 (WrappedeETH): D26, E26
Set-TextValue 'D26' '3.392.28'
Set-TextValue 'E26' '  -3.33%  '

# Row 27 (Litecoin): D27, E27
Set-TextValue 'D27' '74.20'
Set-TextValue 'E27' '  -5.54%  '

# Row 28 (PEPE): E28
Set-TextValue 'E28' '  +2.88%  '

# Row 29 (Dai): E29
Set-TextValue 'E29' '  +0.00%  '

# Row 30 (Cronos): D30, E30
Set-TextValue 'D30' '0.172'
Set-TextValue 'E30' '  -14.57%  '

# Row 31 (Binance-PegBSC-USD): E31
Set-TextValue 'E31' '  +0.16%  '

# Row 32 (Bittensor): D32, E32
Set-TextValue 'D32' '544.48'
Set-TextValue 'E32' '  -10.43%  '

# Row 33 (InternetComputer(DFINITY)): D33, E33
Set-TextValue 'D33' '8.41'
Set-TextValue 'E33' '  -10.84%  '

# Row 34 (PancakeSwap): E34
Set-TextValue 'E34' '  -11.89%  '

# Row 35 (Fetch.AI): E35
Set-TextValue 'E35' '  -20.83%  '

# Row 36 (RenderToken): E36
Set-TextValue 'E36' '  -9.02%  '

# Row 37 (Kaspa): D37, E37
Set-TextValue 'D37' '0.135'
Set-TextValue 'E37' '  -7.75%  '

# Row 38 (EthereumClassic): D38, E38
Set-TextValue 'D38' '22.19'
Set-TextValue 'E38' '  -4.54%  '

# Row 39 (WhiteBITCoin): D39, E39
Set-TextValue 'D39' '21.85'
Set-TextValue 'E39' '  +0.02%  '

# Row 40 (FirstDigitalUSD): D40, E40
Set-TextValue 'D40' '1.00'
Set-TextValue 'E40' '  +0.27%  '

# Row 41 (dogwifhat): D41, E41
Set-TextValue 'D41' '3.00'
Set-TextValue 'E41' '  -2.70%  '

# Row 42 (PolygonEcosystemToken): D42, E42
Set-TextValue 'D42' '0.381'
Set-TextValue 'E42' '  -9.65%  '

# Row 43 (USDe): E43
Set-TextValue 'E43' '  -0.07%  '

# Row 44 (Stacks): D44, E44
Set-TextValue 'D44' '1.89'
Set-TextValue 'E44' '  -14.02%  '

# Row 45 (Monero): D45, E45
Set-TextValue 'D45' '146.37'
Set-TextValue 'E45' '  -8.60%  '

# Row 46 (Aave): D46, E46
Set-TextValue 'D46' '175.10'
Set-TextValue 'E46' '  -8.47%  '

# Row 47 (OKB): B47, C47, D47, E47
Set-TextValue 'B47' 'Stellar'
Set-TextValue 'C47' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D47' '0.129'
Set-TextValue 'E47' '  +12.69%  '

# Row 48 (Stellar): B48, C48, D48, E48
Set-TextValue 'B48' 'OKB'
Set-TextValue 'C48' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D48' '43.48'
Set-TextValue 'E48' '  -6.09%  '

# Row 49 (ImmutableX): D49, E49
Set-TextValue 'D49' '1.27'
Set-TextValue 'E49' '  -11.71%  '

# Row 50 (Filecoin): E50
Set-TextValue 'E50' '  -9.02%  '

# Row 51 (ARBITRUM): D51, E51
Set-TextValue 'D51' '0.601'
Set-TextValue 'E51' '  -9.37%  '
